# Update forecast figures on the "Forecast Comparison" sheet (Auto Arima
# removed -> Amazon Mean Forecast column now populated, other forecast
# columns recalculated) and refresh the dependent totals on "Summary".

$wb  = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New values for columns C (Prophet Forecast), D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90 Forecast)
# for rows 2-17 (weeks W01-W16).
$data = @(
    @(112,80,95,108,129),
    @(102,58,70,82,100),
    @(91,54,64,75,92),
    @(85,55,66,77,95),
    @(85,56,68,81,102),
    @(90,56,68,81,102),
    @(90,57,69,83,106),
    @(88,58,70,85,109),
    @(83,55,66,79,100),
    @(77,55,67,81,103),
    @(70,57,69,84,108),
    @(66,59,72,88,113),
    @(59,56,68,82,105),
    @(49,54,66,81,107),
    @(38,53,65,80,104),
    @(33,51,62,76,100)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $wsForecast.Range("C$row").Value = $vals[0]
    $wsForecast.Range("D$row").Value = $vals[1]
    $wsForecast.Range("E$row").Value = $vals[2]
    $wsForecast.Range("F$row").Value = $vals[3]
    $wsForecast.Range("G$row").Value = $vals[4]
}

# Refresh the rolled-up totals on the Summary sheet to match the new
# forecast figures above. These cells are stored as text, so force a
# text number format to preserve that representation.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1218"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "743"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "390"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "112"
